$wb = $excel.ActiveWorkbook

# --- Markers sheet: add two new drug rows (Hydroxychloroquine) ---
$ws = $wb.Worksheets.Item("Markers")

# Populate column B first, then column A, so new shared strings are
# interned in the same order as the target workbook.
$ws.Range("B8").Value = "Jiný lék na Covid/plaquenil"
$ws.Range("A8").Value = "Hydroxychloroquine"
$ws.Range("B9").Value = "Jiný lék na Covid-plaquenil"
$ws.Range("A9").Value = "Hydroxychloroquine"

# Move the active selection the way it ends up after entering the data.
[void]$ws.Range("B10").Select()

# --- PatientColumns sheet: row height adjustments ---
$ps = $wb.Worksheets.Item("PatientColumns")

$ps.Rows.Item(5).RowHeight = 30
$ps.Rows.Item(6).RowHeight = 30
$ps.Rows.Item(8).RowHeight = 60
$ps.Rows.Item(9).RowHeight = 30
$ps.Rows.Item(10).RowHeight = 45
$ps.Rows.Item(15).RowHeight = 30
$ps.Rows.Item(17).RowHeight = 30
$ps.Rows.Item(18).RowHeight = 30
$ps.Rows.Item(23).AutoFit()
$ps.Rows.Item(24).AutoFit()
$ps.Rows.Item(25).RowHeight = 30
$ps.Rows.Item(26).RowHeight = 45
$ps.Rows.Item(27).RowHeight = 45
$ps.Rows.Item(28).RowHeight = 75
